# This workbook is a weekly price log where each row holds one day's
# price record (columns D, J, K, L, M, N, O, P, Q) for the same
# market/product. The edit described by the diff inserts one new daily
# record at row 4 and shifts every subsequent record (old rows 4-77)
# down by one row (into rows 5-78); the record that used to be in the
# last row (78) now lands in a brand new row 79. Rows 1-3 and the
# descriptive columns (A, B, C, E, F, G, H, I, R), which are identical
# across every data row, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remember the values currently sitting in the last data row (78):
#    after the shift this becomes the new row 79.
$lastRow = $ws.Range("A78:R78").Value2

# 2) Shift the existing rows 4-77 down into rows 5-78 (row 4 itself is
#    left in place by this step since it's only the source of the copy).
$block = $ws.Range("A4:R77").Value2
$ws.Range("A5:R78").Value2 = $block

# 3) Drop the old last row (78) into the newly created row 79, matching
#    the date column's number format used by the rest of column D.
$ws.Range("A79:R79").Value2 = $lastRow
$ws.Range("D79").NumberFormat = $ws.Range("D2").NumberFormat

# 4) Overwrite row 4 with the brand-new daily record. The unit column
#    (N) and the Kg/Unit column (Q) keep their previous values (60),
#    only the date, volume, prices, origin region and $/Kg change.
$ws.Range("D4").Value = 44530
$ws.Range("J4").Value = 350
$ws.Range("K4").Value = 5000
$ws.Range("L4").Value = 5500
$ws.Range("M4").Value = 5214
$ws.Range("O4").Value = "Región del Maule"
$ws.Range("P4").Value = 87
